$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.917.21'
$ws.Range("E2").Value = '  +0.46%  '
$ws.Range("D3").Value = '2.638.81'
$ws.Range("E3").Value = '  +1.27%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.67'
$ws.Range("E5").Value = '  +1.26%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '157.25'
$ws.Range("E6").Value = '  +1.86%  '
$ws.Range("E7").Value = '  +1.06%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.119'
$ws.Range("E9").Value = '  -1.09%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.84'
$ws.Range("E10").Value = '  +1.17%  '
$ws.Range("E11").Value = '  +0.59%  '
$ws.Range("E12").Value = '  +1.16%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '28.84'
$ws.Range("E13").Value = '  +2.56%  '
$ws.Range("D14").Value = '3.112.95'
$ws.Range("E14").Value = '  +1.32%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000185'
$ws.Range("E15").Value = '  +0.95%  '
$ws.Range("D16").Value = '63.831.49'
$ws.Range("E16").Value = '  +0.60%  '
$ws.Range("D17").Value = '2.638.70'
$ws.Range("E17").Value = '  +1.89%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.20'
$ws.Range("E18").Value = '  +1.33%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.77'
$ws.Range("E19").Value = '  +3.89%  '
$ws.Range("E20").Value = '  -1.80%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '344.42'
$ws.Range("E22").Value = '  +0.18%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.38'
$ws.Range("E23").Value = '  +2.22%  '
$ws.Range("E24").Value = '  +9.52%  '
$ws.Range("E25").Value = '  +5.31%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.66'
$ws.Range("E26").Value = '  +5.12%  '
$ws.Range("E27").Value = '  +0.62%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '585.33'
$ws.Range("E28").Value = '  +1.71%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.24'
$ws.Range("E29").Value = '  +4.83%  '
$ws.Range("E30").Value = '  +0.85%  '
$ws.Range("E31").Value = '  +0.13%  '
$ws.Range("E32").Value = '  -0.12%  '
$ws.Range("E33").Value = '  +2.83%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.65'
$ws.Range("E34").Value = '  +3.43%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.48'
$ws.Range("E35").Value = '  +3.23%  '
$ws.Range("E36").Value = '  -0.70%  '
$ws.Range("E37").Value = '  +0.00%  '
$ws.Range("E38").Value = '  +0.00%  '
$ws.Range("E39").Value = '  +3.57%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '153.71'
$ws.Range("E40").Value = '  +1.11%  '
$ws.Range("E41").Value = '  +10.07%  '
$ws.Range("E42").Value = '  +0.06%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '162.58'
$ws.Range("E43").Value = '  +4.38%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '24.22'
$ws.Range("E44").Value = '  +6.31%  '
$ws.Range("E46").Value = '  -0.10%  '
$ws.Range("E47").Value = '  +1.26%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.101'
$ws.Range("E48").Value = '  -0.50%  '
$ws.Range("E49").Value = '  +0.17%  '
$ws.Range("D50").Value = '0.0₆0238'
$ws.Range("E50").Value = '  +1.99%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.795'
$ws.Range("E51").Value = '  +3.26%  '
